$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.872739
$ws.Range("H2").Value = 7.745478
$ws.Range("I2").Value = 0.03299520440786341
$ws.Range("J2").Value = 0.02275854185403964
$ws.Range("Q2").Value = 1.079128395046
$ws.Range("R2").Value = 6.474770370276
$ws.Range("S2").Value = 0.03299520440786341
$ws.Range("T2").Value = 0.02275854185403964

$ws.Range("I3").Value = 0.8633805772213771
$ws.Range("J3").Value = 0.8932790395734661
$ws.Range("S3").Value = 0.8633805772213771
$ws.Range("T3").Value = 0.8932790395734661

$ws.Range("G4").Value = 2.422325
$ws.Range("H4").Value = 7.266975
$ws.Range("I4").Value = 0.02063787632403778
$ws.Range("J4").Value = 0.02135255625150052
$ws.Range("Q4").Value = 0.6749744017166668
$ws.Range("R4").Value = 6.074769615450001
$ws.Range("S4").Value = 0.02063787632403778
$ws.Range("T4").Value = 0.02135255625150052

$ws.Range("G5").Value = 7.912825
$ws.Range("H5").Value = 15.82565
$ws.Range("I5").Value = 0.06741618227271494
$ws.Range("J5").Value = 0.04650051525449849
$ws.Range("Q5").Value = 2.204887585383334
$ws.Range("R5").Value = 13.2293255123
$ws.Range("S5").Value = 0.06741618227271494
$ws.Range("T5").Value = 0.04650051525449849

$ws.Range("G6").Value = 0.5045936666666666
$ws.Range("H6").Value = 1.513781
$ws.Range("I6").Value = 0.004299068740387607
$ws.Range("J6").Value = 0.004447943464089625
$ws.Range("Q6").Value = 0.1406036796335556
$ws.Range("R6").Value = 1.265433116702
$ws.Range("S6").Value = 0.004299068740387607
$ws.Range("T6").Value = 0.004447943464089625

$ws.Range("G7").Value = 1.322919333333333
$ws.Range("H7").Value = 3.968758
$ws.Range("I7").Value = 0.01127109103361929
$ws.Range("J7").Value = 0.01166140360240577
$ws.Range("Q7").Value = 0.3686279444484445
$ws.Range("R7").Value = 3.317651500036
$ws.Range("S7").Value = 0.01127109103361929
$ws.Range("T7").Value = 0.01166140360240577
